# ISB_Forecasting_Tutorials.pptx edit script
# 1) Update the footer "datetimeFigureOut" field placeholders (slide layouts +
#    slide master) from 1/31/2018 to 2/4/2018.
# 2) Replace the "02/04/2018" line on slide 1 with an "Email: ..." paragraph
#    and a "Code: ..." paragraph, each containing a hyperlinked run.

$p = $ppt.ActivePresentation

function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "1/31/2018") {
                $tr.Text = "2/4/2018"
            }
        }
    }
}

# --- Slide master ---
Update-DateField $p.SlideMaster.Shapes

# --- Slide layouts ---
for ($li = 1; $li -le $p.SlideMaster.CustomLayouts.Count; $li++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

# --- Slide 1: swap the date line for Email/Code contact paragraphs ---
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$firstLine = "Suryanarayana, Data Scientist, Flex "
$emailLabel = "Email: "
$email = "Suryanarayana_Ambatipudi_2014@cba.isb.edu"
$codeLabel = "Code: "
$codeUrl = "https://github.com/asuryam/ISB-Tutorials/tree/Forecasting"
$trailingSpace = " "

$tr.Text = $firstLine + "`r" + $emailLabel + $email + "`r" + $codeLabel + $codeUrl + $trailingSpace

# Compute absolute 1-based character offsets by hand (the runtime's
# Paragraphs()/Runs() sub-ranges don't report Start/Length here), each
# paragraph break ("`r") also counts as one character.
$emailStart = $firstLine.Length + 1 + $emailLabel.Length + 1
$emailRange = $tr.Characters($emailStart, $email.Length)
$emailRange.ActionSettings(1).Hyperlink.Address = "mailto:" + $email

$codeStart = $emailStart + $email.Length + 1 + $codeLabel.Length
$codeRange = $tr.Characters($codeStart, $codeUrl.Length)
$codeRange.ActionSettings(1).Hyperlink.Address = $codeUrl
